# Commit: "added cluster setup steps"
# Change the SASRec (PyTorch) results bullet on slide 2 from
# " - As expected." to " - As in paper.", keeping the surrounding
# run (formatting) untouched.

$p = $ppt.ActivePresentation

$oldFragment = "As expected."
$newFragment = "As in paper."

# The run that needs editing starts right after "PyTorch" with the
# en-dash " - Legit?): Successfully ..." and ends with the closing
# "... numerical stability?)". Using those unique anchors lets us
# grab the *entire* run as a single Characters() range, so setting
# its .Text keeps it as one <a:r> instead of splintering the run.
$runStartAnchor = " – Legit?): Successfully"
$runEndAnchor   = "numerical stability?)"

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($shi = 1; $shi -le $slide.Shapes.Count; $shi++) {
        $shape = $slide.Shapes.Item($shi)
        if (-not $shape.HasTextFrame) { continue }

        $tr = $shape.TextFrame.TextRange
        $full = $tr.Text

        if ($full.IndexOf($oldFragment) -lt 0) { continue }
        if ($full.IndexOf($runStartAnchor) -lt 0) { continue }

        $startIdx = $full.IndexOf($runStartAnchor)
        $endAnchorIdx = $full.IndexOf($runEndAnchor)
        $runLen = ($endAnchorIdx + $runEndAnchor.Length) - $startIdx

        $runRange = $tr.Characters($startIdx + 1, $runLen)
        $runText = $runRange.Text
        $runRange.Text = $runText.Replace($oldFragment, $newFragment)
    }
}
